$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "38.856.94"
$ws.Cells.Item(2, 5).Value = "  +0.14%  "

# Row 3
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = "2.140.06"
$ws.Cells.Item(3, 5).Value = "  +2.31%  "

# Row 4
$ws.Cells.Item(4, 5).Value = "  +0.12%  "

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "228.27"
$ws.Cells.Item(5, 5).Value = "  -0.06%  "

# Row 6
$ws.Cells.Item(6, 5).Value = "  +0.53%  "

# Row 7
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "62.24"
$ws.Cells.Item(7, 5).Value = "  +2.36%  "

# Row 8
$ws.Cells.Item(8, 5).Value = "  +0.01%  "

# Row 9
$ws.Cells.Item(9, 5).Value = "  +1.48%  "

# Row 10
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "0.0845"
$ws.Cells.Item(10, 5).Value = "  +1.00%  "

# Row 11
$ws.Cells.Item(11, 5).Value = "  -0.57%  "

# Row 12
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "15.97"
$ws.Cells.Item(12, 5).Value = "  +6.58%  "

# Row 13
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "2.445.52"
$ws.Cells.Item(13, 5).Value = "  +1.79%  "

# Row 14
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "22.14"
$ws.Cells.Item(14, 5).Value = "  +0.81%  "

# Row 15
$ws.Cells.Item(15, 5).Value = "  +1.36%  "

# Row 16
$ws.Cells.Item(16, 5).Value = "  +0.66%  "

# Row 17
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "2.155.83"
$ws.Cells.Item(17, 5).Value = "  +3.23%  "

# Row 18
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "38.885.36"
$ws.Cells.Item(18, 5).Value = "  +0.36%  "

# Row 19
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "71.86"
$ws.Cells.Item(19, 5).Value = "  +0.03%  "

# Row 20
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "6.12"
$ws.Cells.Item(20, 5).Value = "  +1.71%  "

# Row 21
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "0.0₃0847"
$ws.Cells.Item(21, 5).Value = "  +1.33%  "

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "227.68"
$ws.Cells.Item(22, 5).Value = "  +0.54%  "

# Row 23
$ws.Cells.Item(23, 5).Value = "  -0.06%  "

# Row 24
$ws.Cells.Item(24, 5).Value = "  -3.98%  "

# Row 25
$ws.Cells.Item(25, 5).Value = "  -0.58%  "

# Row 26
$ws.Cells.Item(26, 5).Value = "  +2.63%  "

# Row 27
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "170.73"
$ws.Cells.Item(27, 5).Value = "  +0.04%  "

# Row 28
$ws.Cells.Item(28, 5).Value = "  +0.21%  "

# Row 29
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "1.41"
$ws.Cells.Item(29, 5).Value = "  -2.88%  "

# Row 30
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "19.53"
$ws.Cells.Item(30, 5).Value = "  +1.88%  "

# Row 31
$ws.Cells.Item(31, 5).Value = "  +9.00%  "

# Row 32
$ws.Cells.Item(32, 5).Value = "  +0.71%  "

# Row 33
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "4.60"
$ws.Cells.Item(33, 5).Value = "  +2.48%  "

# Row 34
$ws.Cells.Item(34, 2).Value = "InternetComputer(DFINITY)"
$ws.Cells.Item(34, 3).Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "4.80"
$ws.Cells.Item(34, 5).Value = "  +2.14%  "

# Row 35
$ws.Cells.Item(35, 2).Value = "THORChain"
$ws.Cells.Item(35, 3).Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "7.17"
$ws.Cells.Item(35, 5).Value = "  +11.43%  "

# Row 36
$ws.Cells.Item(36, 5).Value = "  +0.62%  "

# Row 37
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "2.40"
$ws.Cells.Item(37, 5).Value = "  +0.13%  "

# Row 38
$ws.Cells.Item(38, 5).Value = "  +0.42%  "

# Row 39
$ws.Cells.Item(39, 5).Value = "  +0.10%  "

# Row 40
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "18.16"
$ws.Cells.Item(40, 5).Value = "  -0.42%  "

# Row 41
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "0.0230"
$ws.Cells.Item(41, 5).Value = "  +3.27%  "

# Row 42
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "102.70"
$ws.Cells.Item(42, 5).Value = "  +1.57%  "

# Row 43
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "1.532.04"
$ws.Cells.Item(43, 5).Value = "  -0.45%  "

# Row 44
$ws.Cells.Item(44, 5).Value = "  +6.68%  "

# Row 45
$ws.Cells.Item(45, 5).Value = "  +2.49%  "

# Row 46
$ws.Cells.Item(46, 5).Value = "  -0.58%  "

# Row 47
$ws.Cells.Item(47, 5).Value = "  +5.89%  "

# Row 48
$ws.Cells.Item(48, 5).Value = "  -0.96%  "

# Row 49
$ws.Cells.Item(49, 5).Value = "  +1.44%  "

# Row 50
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "2.98"
$ws.Cells.Item(50, 5).Value = "  +0.21%  "

# Row 51
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "2.330.18"
$ws.Cells.Item(51, 5).Value = "  +1.82%  "
